# Auto-generated Excel COM-interop script applying the diff's cell value changes.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 772.5
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 772.5
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 2317.5
$ws.Range("N17").Value = -2653.5
$ws.Range("H32").Value = 5795.2
$ws.Range("I32").Value = 7568.2
$ws.Range("J32").Value = 4908.7
$ws.Range("K32").Value = 7568.2
$ws.Range("L32").Value = 4908.7
$ws.Range("M32").Value = -7242.2
$ws.Range("N32").Value = -5560.7
$ws.Range("H40").Value = 3199.125
$ws.Range("I40").Value = 2320
$ws.Range("J40").Value = 4664.3335
$ws.Range("K40").Value = 2320
$ws.Range("L40").Value = 4664.3335
$ws.Range("M40").Value = -2145
$ws.Range("N40").Value = -5014.3335
$ws.Range("H53").Value = 467.7
$ws.Range("I53").Value = 459.66666
$ws.Range("J53").Value = 471.14285
$ws.Range("K53").Value = 459.66666
$ws.Range("L53").Value = 471.14285
$ws.Range("M53").Value = 177.33334
$ws.Range("N53").Value = -1745.14285
$ws.Range("H106").Value = 1689.6923
$ws.Range("I106").Value = 1497.1666
$ws.Range("J106").Value = 4000
$ws.Range("K106").Value = 1497.1666
$ws.Range("L106").Value = 4000
$ws.Range("M106").Value = -866.1666
$ws.Range("N106").Value = -5262
$ws.Range("H111").Value = 200
$ws.Range("I111").Value = 200
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 600
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = 2467
$ws.Range("H132").Value = 18117.572
$ws.Range("I132").Value = 19139.46
$ws.Range("J132").Value = 4833
$ws.Range("K132").Value = 57418.38
$ws.Range("L132").Value = 14499
$ws.Range("M132").Value = -54888.38
$ws.Range("N132").Value = -19559
$ws.Range("H134").Value = 99999
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 99999
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 99999
$ws.Range("N134").Value = -110139
$ws.Range("H138").Value = 31989.383
$ws.Range("I138").Value = 2502.9443
$ws.Range("J138").Value = 65161.625
$ws.Range("K138").Value = 7508.8329
$ws.Range("L138").Value = 195484.875
$ws.Range("M138").Value = -2368.8329
$ws.Range("N138").Value = -205764.875
$ws.Range("H141").Value = 3171.4285
$ws.Range("I141").Value = 3400
$ws.Range("J141").Value = 1800
$ws.Range("K141").Value = 10200
$ws.Range("L141").Value = 5400
$ws.Range("M141").Value = -5020
$ws.Range("N141").Value = -15760

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 33235.438
$ws.Range("I32").Value = 35351.133
$ws.Range("J32").Value = 1500
$ws.Range("K32").Value = 35351.133
$ws.Range("L32").Value = 1500
$ws.Range("M32").Value = -35064.133
$ws.Range("N32").Value = -2074
$ws.Range("H110").Value = 593.63635
$ws.Range("I110").Value = 553
$ws.Range("J110").Value = 1000
$ws.Range("K110").Value = 553
$ws.Range("L110").Value = 1000
$ws.Range("M110").Value = 1492
$ws.Range("N110").Value = -5090
$ws.Range("H132").Value = 2272
$ws.Range("I132").Value = 1998.8572
$ws.Range("J132").Value = 2750
$ws.Range("K132").Value = 5996.571599999999
$ws.Range("L132").Value = 8250
$ws.Range("M132").Value = -3466.571599999999
$ws.Range("N132").Value = -13310

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H138").Value = 40000
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 40000
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 40000
$ws.Range("N138").Value = -50280
$ws.Range("H139").Value = 111625
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 111625
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 111625
$ws.Range("N139").Value = -121905

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5264359.5
$ws.Range("I31").Value = 7693225
$ws.Range("J31").Value = 1816.6666
$ws.Range("K31").Value = 7693225
$ws.Range("L31").Value = 1816.6666
$ws.Range("M31").Value = -7692930
$ws.Range("N31").Value = -2406.6666
$ws.Range("H34").Value = 5264359.5
$ws.Range("I34").Value = 7693225
$ws.Range("J34").Value = 1816.6666
$ws.Range("K34").Value = 7693225
$ws.Range("L34").Value = 1816.6666
$ws.Range("M34").Value = -7693023
$ws.Range("N34").Value = -2220.6666
$ws.Range("H58").Value = 1601.6666
$ws.Range("I58").Value = 1191.3636
$ws.Range("J58").Value = 2730
$ws.Range("K58").Value = 1191.3636
$ws.Range("L58").Value = 2730
$ws.Range("M58").Value = -988.3635999999999
$ws.Range("N58").Value = -3136
$ws.Range("H99").Value = 8998.25
$ws.Range("I99").Value = 5748.125
$ws.Range("J99").Value = 15498.5
$ws.Range("K99").Value = 5748.125
$ws.Range("L99").Value = 15498.5
$ws.Range("M99").Value = -4250.125
$ws.Range("N99").Value = -18494.5
$ws.Range("H122").Value = 1851.3077
$ws.Range("I122").Value = 1833.875
$ws.Range("J122").Value = 1879.2
$ws.Range("K122").Value = 5501.625
$ws.Range("L122").Value = 5637.6
$ws.Range("M122").Value = -3051.625
$ws.Range("N122").Value = -10537.6
$ws.Range("H126").Value = 8998.25
$ws.Range("I126").Value = 5748.125
$ws.Range("J126").Value = 15498.5
$ws.Range("K126").Value = 17244.375
$ws.Range("L126").Value = 46495.5
$ws.Range("M126").Value = -14774.375
$ws.Range("N126").Value = -51435.5
$ws.Range("H127").Value = 54999.668
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 54999.668
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 54999.668
$ws.Range("N127").Value = -64919.668
$ws.Range("H132").Value = 60331.65
$ws.Range("I132").Value = 84262
$ws.Range("J132").Value = 2898.8
$ws.Range("K132").Value = 252786
$ws.Range("L132").Value = 8696.400000000001
$ws.Range("M132").Value = -250256
$ws.Range("N132").Value = -13756.4
$ws.Range("H135").Value = 119994
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 119994
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 119994
$ws.Range("N135").Value = -130134
$ws.Range("H136").Value = 1601.6666
$ws.Range("I136").Value = 1191.3636
$ws.Range("J136").Value = 2730
$ws.Range("K136").Value = 3574.0908
$ws.Range("L136").Value = 8190
$ws.Range("M136").Value = -1024.0908
$ws.Range("N136").Value = -13290
$ws.Range("H140").Value = 105877.766
$ws.Range("I140").Value = 99999
$ws.Range("J140").Value = 108327.25
$ws.Range("K140").Value = 99999
$ws.Range("L140").Value = 108327.25
$ws.Range("M140").Value = -94819
$ws.Range("N140").Value = -118687.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 1106.6666
$ws.Range("I114").Value = 125
$ws.Range("J114").Value = 1597.5
$ws.Range("K114").Value = 375
$ws.Range("L114").Value = 4792.5
$ws.Range("M114").Value = 2879
$ws.Range("N114").Value = -11300.5
$ws.Range("H129").Value = 3155.4736
$ws.Range("I129").Value = 2108.3333
$ws.Range("J129").Value = 4950.5713
$ws.Range("K129").Value = 6324.999899999999
$ws.Range("L129").Value = 14851.7139
$ws.Range("M129").Value = -1324.999899999999
$ws.Range("N129").Value = -24851.7139

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 100004600
$ws.Range("I122").Value = 5500
$ws.Range("J122").Value = 250003250
$ws.Range("K122").Value = 16500
$ws.Range("L122").Value = 750009750
$ws.Range("M122").Value = -14050
$ws.Range("N122").Value = -750014650

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4541.9
$ws.Range("I7").Value = 3713.7144
$ws.Range("J7").Value = 6474.3335
$ws.Range("K7").Value = 3713.7144
$ws.Range("L7").Value = 6474.3335
$ws.Range("M7").Value = -3601.7144
$ws.Range("N7").Value = -6698.3335
$ws.Range("H61").Value = 3401.25
$ws.Range("I61").Value = 3185
$ws.Range("J61").Value = 4050
$ws.Range("K61").Value = 3185
$ws.Range("L61").Value = 4050
$ws.Range("M61").Value = -2983
$ws.Range("N61").Value = -4454
$ws.Range("H68").Value = 4738.722
$ws.Range("I68").Value = 4574.625
$ws.Range("J68").Value = 4870
$ws.Range("K68").Value = 4574.625
$ws.Range("L68").Value = 4870
$ws.Range("M68").Value = -3825.625
$ws.Range("N68").Value = -6368
$ws.Range("H71").Value = 4738.722
$ws.Range("I71").Value = 4574.625
$ws.Range("J71").Value = 4870
$ws.Range("K71").Value = 22873.125
$ws.Range("L71").Value = 24350
$ws.Range("M71").Value = -19129.125
$ws.Range("N71").Value = -31838
$ws.Range("H113").Value = 3401.25
$ws.Range("I113").Value = 3185
$ws.Range("J113").Value = 4050
$ws.Range("K113").Value = 3185
$ws.Range("L113").Value = 4050
$ws.Range("M113").Value = -1015
$ws.Range("N113").Value = -8390
$ws.Range("H122").Value = 3955.4443
$ws.Range("I122").Value = 3957
$ws.Range("J122").Value = 3950
$ws.Range("K122").Value = 11871
$ws.Range("L122").Value = 11850
$ws.Range("M122").Value = -9421
$ws.Range("N122").Value = -16750
$ws.Range("H126").Value = 4541.9
$ws.Range("I126").Value = 3713.7144
$ws.Range("J126").Value = 6474.3335
$ws.Range("K126").Value = 11141.1432
$ws.Range("L126").Value = 19423.0005
$ws.Range("M126").Value = -8671.143199999999
$ws.Range("N126").Value = -24363.0005
$ws.Range("H132").Value = 4813.6
$ws.Range("I132").Value = 4564.067
$ws.Range("J132").Value = 5187.9
$ws.Range("K132").Value = 13692.201
$ws.Range("L132").Value = 15563.7
$ws.Range("M132").Value = -11162.201
$ws.Range("N132").Value = -20623.7
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("M137").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4721.5
$ws.Range("I62").Value = 4444
$ws.Range("J62").Value = 4999
$ws.Range("K62").Value = 4444
$ws.Range("L62").Value = 4999
$ws.Range("M62").Value = -3820
$ws.Range("N62").Value = -6247
$ws.Range("H65").Value = 4721.5
$ws.Range("I65").Value = 4444
$ws.Range("J65").Value = 4999
$ws.Range("K65").Value = 22220
$ws.Range("L65").Value = 24995
$ws.Range("M65").Value = -19100
$ws.Range("N65").Value = -31235
$ws.Range("H107").Value = 1237.2307
$ws.Range("I107").Value = 1159
$ws.Range("J107").Value = 1498
$ws.Range("K107").Value = 3477
$ws.Range("L107").Value = 4494
$ws.Range("M107").Value = -1557
$ws.Range("N107").Value = -8334
$ws.Range("H119").Value = 45000
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 45000
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 45000
$ws.Range("N119").Value = -54676
$ws.Range("H122").Value = 121431
$ws.Range("I122").Value = 132125.19
$ws.Range("J122").Value = 3795
$ws.Range("K122").Value = 396375.57
$ws.Range("L122").Value = 11385
$ws.Range("M122").Value = -393925.57
$ws.Range("N122").Value = -16285
$ws.Range("H126").Value = 458577.9
$ws.Range("I126").Value = 2891
$ws.Range("J126").Value = 1005402.2
$ws.Range("K126").Value = 8673
$ws.Range("L126").Value = 3016206.6
$ws.Range("M126").Value = -6203
$ws.Range("N126").Value = -3021146.6
$ws.Range("H132").Value = 24411.219
$ws.Range("I132").Value = 37008.2
$ws.Range("J132").Value = 3416.25
$ws.Range("K132").Value = 111024.6
$ws.Range("L132").Value = 10248.75
$ws.Range("M132").Value = -108494.6
$ws.Range("N132").Value = -15308.75

Write-Host "Applied all cell value updates."